$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows for 15, 16, 17 March 2021 (serial dates 44301-44303)
$data = @(
    @(44301, 0, 5, 187.3360809291869),
    @(44302, 0, 4, 149.8688647433496),
    @(44303, 0, 3, 112.4016485575122)
)

$startRow = 227

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    # Copy the formatting of column A from the previous row (date style) so the
    # new date cell matches the existing look (border/bold/center/date format)
    $ws.Range("A" + ($row - 1)).Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$excel.CutCopyMode = 0
